$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D9").Value = "“강의 어렵고 학점 짜다” 뉴욕대, 학생들 불만에 교수 잘랐다"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/nyu-shame/#utm_source=rss&utm_medium=rss&utm_campaign=nyu-shame"

$ws.Range("D28").Value = "Mobile Manipulator 101 :: Pseudo-Inverse and Singularity-Robust Inverse"
$ws.Range("E28").Value = "https://ropiens.tistory.com/198"

$ws.Range("D51").Value = "영어 표현 more than은 기준 값을 포함하지 않는다 (게임회사 K사 코테 준비)"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EC%98%81%EC%96%B4-%ED%91%9C%ED%98%84-more-than%EC%9D%80-%EA%B8%B0%EC%A4%80-%EA%B0%92%EC%9D%84-%ED%8F%AC%ED%95%A8%ED%95%98%EC%A7%80-%EC%95%8A%EB%8A%94%EB%8B%A4-%EA%B2%8C%EC%9E%84%ED%9A%8C%EC%82%AC-K%EC%82%AC-%EC%BD%94%ED%85%8C-%EC%A4%80%EB%B9%84"
